$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; C=324074; D=412922472},
    @{Row=3; C=261; D=311952},
    @{Row=4; C=324; D=464207},
    @{Row=8; C=869; D=1278295},
    @{Row=10; C=117874; D=172715765},
    @{Row=12; C=60069; D=86696408},
    @{Row=16; C=4025; D=5712197},
    @{Row=20; C=6830; D=9533037},
    @{Row=22; C=78190; D=97445489},
    @{Row=28; C=32667; D=47814776},
    @{Row=30; C=11577; D=16653312},
    @{Row=33; C=1564; D=2197281},
    @{Row=35; C=1880; D=2655168},
    @{Row=36; C=98050; D=123337544},
    @{Row=44; C=44611; D=65380019},
    @{Row=46; C=9231; D=13243911},
    @{Row=48; C=1414; D=1963603},
    @{Row=51; C=2375; D=3317952},
    @{Row=52; C=69681; D=87396051},
    @{Row=59; C=28380; D=41621428},
    @{Row=62; C=11236; D=16247212},
    @{Row=64; C=1365; D=1907737},
    @{Row=68; C=1514; D=2122704},
    @{Row=70; C=20678; D=27080493},
    @{Row=74; C=7655; D=11210276},
    @{Row=76; C=5168; D=7504063},
    @{Row=78; C=284; D=399083},
    @{Row=79; C=142179; D=177218222},
    @{Row=83; C=431; D=629324},
    @{Row=85; C=63966; D=93751771},
    @{Row=88; C=30008; D=43406845},
    @{Row=90; C=2746; D=3953457},
    @{Row=91; C=2893; D=4090360},
    @{Row=92; C=33679; D=45646991},
    @{Row=96; C=8158; D=11994580},
    @{Row=98; C=7509; D=10893587},
    @{Row=101; C=501; D=722891},
    @{Row=102; C=10434; D=15915355},
    @{Row=104; C=2573; D=4195321},
    @{Row=106; C=3478; D=5664690},
    @{Row=110; C=142860; D=176655196},
    @{Row=116; C=53104; D=77837043},
    @{Row=118; C=27413; D=39719852},
    @{Row=119; C=1317; D=1800981},
    @{Row=122; C=2317; D=3257168},
    @{Row=124; C=521312; D=688532394},
    @{Row=126; C=218; D=321009},
    @{Row=129; C=1386; D=2054182},
    @{Row=131; C=210405; D=309309006},
    @{Row=134; C=185957; D=270422504},
    @{Row=137; C=2863; D=4021893},
    @{Row=139; C=6517; D=9208214},
    @{Row=142; C=45222; D=60373149},
    @{Row=148; C=14209; D=20831928},
    @{Row=149; C=3804; D=5485240},
    @{Row=152; C=404; D=581216},
    @{Row=154; C=395; D=558263},
    @{Row=155; C=17825; D=23559950},
    @{Row=159; C=7275; D=10585075},
    @{Row=161; C=5074; D=7303496},
    @{Row=166; C=18533; D=30115202},
    @{Row=167; C=2030; D=3326572},
    @{Row=172; C=88621; D=110760198},
    @{Row=179; C=34116; D=50026949},
    @{Row=181; C=13168; D=19024845},
    @{Row=185; C=1692; D=2378192},
    @{Row=187; C=240589; D=298985378},
    @{Row=189; C=169; D=243736},
    @{Row=195; C=87174; D=127783392},
    @{Row=198; C=33330; D=47979310},
    @{Row=201; C=5150; D=7333077},
    @{Row=204; C=5002; D=6927484},
    @{Row=207; C=266535; D=329806749},
    @{Row=208; C=160; D=176413},
    @{Row=209; C=256; D=365539},
    @{Row=214; C=620; D=902878},
    @{Row=216; C=95808; D=140163257},
    @{Row=217; C=91; D=135699},
    @{Row=219; C=51981; D=75127469},
    @{Row=222; C=4702; D=6600691},
    @{Row=225; C=5916; D=8189397},
    @{Row=228; C=107378; D=134233300},
    @{Row=235; C=49860; D=73039747},
    @{Row=237; C=12582; D=18089606},
    @{Row=239; C=1901; D=2725382},
    @{Row=241; C=2576; D=3604686},
    @{Row=242; C=260526; D=328954344},
    @{Row=248; C=836; D=1227404},
    @{Row=250; C=96556; D=141480229},
    @{Row=253; C=65842; D=95419882},
    @{Row=255; C=2426; D=3420273},
    @{Row=258; C=4715; D=6626760}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
